$wb = $excel.ActiveWorkbook

# Add data for 2025-09-21
# Updates column L (2025 YTD) violent-crime counts across the Citywide
# Totals sheet, the By Neighborhood summary sheet, and every individual
# neighborhood sheet affected by that day of data. A couple of 2024 (K)
# values also receive a one-off correction per the source diff.

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 4888
$ws.Range('L3').Value = 5255
$ws.Range('K4').Value = 1313
$ws.Range('L4').Value = 1284
$ws.Range('L5').Value = 308
$ws.Range('L6').Value = 4440
$ws.Range('K7').Value = 20265
$ws.Range('L7').Value = 16175

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 315
$ws.Range('L3').Value = 365
$ws.Range('L7').Value = 1077

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L4').Value = 14
$ws.Range('L7').Value = 358

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 202
$ws.Range('L3').Value = 257
$ws.Range('L5').Value = 16
$ws.Range('L7').Value = 742

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L2').Value = 183
$ws.Range('L3').Value = 206
$ws.Range('L7').Value = 611

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L2').Value = 114
$ws.Range('L3').Value = 99
$ws.Range('L4').Value = 16
$ws.Range('L6').Value = 84
$ws.Range('L7').Value = 318

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('L2').Value = 75
$ws.Range('L7').Value = 277

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L7').Value = 529
$ws.Range('L8').Value = 1077
$ws.Range('L18').Value = 114
$ws.Range('L20').Value = 406
$ws.Range('L22').Value = 48
$ws.Range('L23').Value = 176
$ws.Range('L27').Value = 146
$ws.Range('L29').Value = 883
$ws.Range('L32').Value = 22
$ws.Range('L33').Value = 742
$ws.Range('L36').Value = 211
$ws.Range('L37').Value = 611
$ws.Range('L40').Value = 44
$ws.Range('L42').Value = 527
$ws.Range('L43').Value = 119
$ws.Range('L48').Value = 211
$ws.Range('L50').Value = 81
$ws.Range('L51').Value = 205
$ws.Range('L54').Value = 341
$ws.Range('L55').Value = 155
$ws.Range('L58').Value = 10
$ws.Range('K63').Value = 152
$ws.Range('L63').Value = 44
$ws.Range('L65').Value = 318
$ws.Range('L66').Value = 40
$ws.Range('L67').Value = 561
$ws.Range('L73').Value = 126
$ws.Range('L76').Value = 249
$ws.Range('L77').Value = 107
$ws.Range('L79').Value = 427
$ws.Range('L83').Value = 358
$ws.Range('L84').Value = 158
$ws.Range('L85').Value = 828
$ws.Range('L86').Value = 117
$ws.Range('L89').Value = 234
$ws.Range('L92').Value = 47
$ws.Range('L94').Value = 196
$ws.Range('L96').Value = 183
$ws.Range('L99').Value = 277
$ws.Range('K101').Value = 20265
$ws.Range('L101').Value = 16175

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L2').Value = 164
$ws.Range('L3').Value = 216
$ws.Range('L5').Value = 13
$ws.Range('L6').Value = 129
$ws.Range('L7').Value = 561

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('L2').Value = 51
$ws.Range('L7').Value = 158

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L2').Value = 62
$ws.Range('L6').Value = 165
$ws.Range('L7').Value = 341

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L2').Value = 267
$ws.Range('L3').Value = 333
$ws.Range('L6').Value = 227
$ws.Range('L7').Value = 883

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('L4').Value = 43
$ws.Range('L6').Value = 88
$ws.Range('L7').Value = 211

$ws = $wb.Worksheets.Item('River North')
$ws.Range('L3').Value = 46
$ws.Range('L7').Value = 249

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L6').Value = 146
$ws.Range('L7').Value = 527

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('L3').Value = 49
$ws.Range('L6').Value = 45
$ws.Range('L7').Value = 155

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('L2').Value = 44
$ws.Range('L3').Value = 67
$ws.Range('L7').Value = 176

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L3').Value = 54
$ws.Range('L6').Value = 52
$ws.Range('L7').Value = 183

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L2').Value = 139
$ws.Range('L3').Value = 152
$ws.Range('L6').Value = 94
$ws.Range('L7').Value = 427

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L2').Value = 125
$ws.Range('L3').Value = 132
$ws.Range('L6').Value = 109
$ws.Range('L7').Value = 406

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('L2').Value = 41
$ws.Range('L7').Value = 114

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('L3').Value = 63
$ws.Range('L7').Value = 211

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L2').Value = 181
$ws.Range('L5').Value = 11
$ws.Range('L6').Value = 127
$ws.Range('L7').Value = 529

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('L4').Value = 25
$ws.Range('L7').Value = 196

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('L3').Value = 20
$ws.Range('L7').Value = 81

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('L3').Value = 10
$ws.Range('L7').Value = 40

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('L4').Value = 11
$ws.Range('L7').Value = 126

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('L2').Value = 18
$ws.Range('L7').Value = 47

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range('L2').Value = 15
$ws.Range('L7').Value = 22

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('L3').Value = 70
$ws.Range('L6').Value = 63
$ws.Range('L7').Value = 234

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('L6').Value = 42
$ws.Range('L7').Value = 146

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('L3').Value = 23
$ws.Range('L7').Value = 117

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('L2').Value = 61
$ws.Range('L7').Value = 205

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('L3').Value = 38
$ws.Range('L7').Value = 119

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L3').Value = 336
$ws.Range('L6').Value = 176
$ws.Range('L7').Value = 828

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('L2').Value = 18
$ws.Range('L7').Value = 48

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('L3').Value = 37
$ws.Range('L7').Value = 107

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range('L3').Value = 20
$ws.Range('L7').Value = 44

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range('L6').Value = 4
$ws.Range('L7').Value = 10
